$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure all touched cells remain plain text, matching the inline-string
# cell type used throughout this sheet (prevents Excel from auto-coercing
# numeric-looking strings such as "214.74" or "1.00" into numbers).

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '25.959.34'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +0.17%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.637.44'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -0.01%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.13%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '214.74'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -0.12%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.505'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -0.21%  '
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +0.25%  '
$ws.Range('B8').NumberFormat = '@'
$ws.Range('B8').Value = 'Dogecoin'
$ws.Range('C8').NumberFormat = '@'
$ws.Range('C8').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.0635'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -0.45%  '
$ws.Range('B9').NumberFormat = '@'
$ws.Range('B9').Value = 'Cardano'
$ws.Range('C9').NumberFormat = '@'
$ws.Range('C9').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.252'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -1.07%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.42'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -0.84%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0794'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +0.10%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '4.23'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -0.50%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.615.91'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -2.61%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.540'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -0.61%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '62.99'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +0.51%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0₃0756'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -0.88%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '25.989.67'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +0.28%  '
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +0.23%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '192.84'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -0.16%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.29'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -1.80%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '9.72'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -1.76%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.14'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -1.97%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.131'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +3.41%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '143.72'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -0.21%  '
$ws.Range('B25').NumberFormat = '@'
$ws.Range('B25').Value = 'Toncoin'
$ws.Range('C25').NumberFormat = '@'
$ws.Range('C25').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.77'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -0.16%  '
$ws.Range('B26').NumberFormat = '@'
$ws.Range('B26').Value = 'BinanceUSD'
$ws.Range('C26').NumberFormat = '@'
$ws.Range('C26').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.00'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +0.08%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '6.82'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -0.16%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '15.43'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -0.34%  '
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -0.16%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0489'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -2.39%  '
$ws.Range('B31').NumberFormat = '@'
$ws.Range('B31').Value = 'Filecoin'
$ws.Range('C31').NumberFormat = '@'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.23'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +0.28%  '
$ws.Range('B32').NumberFormat = '@'
$ws.Range('B32').Value = 'InternetComputer(DFINITY)'
$ws.Range('C32').NumberFormat = '@'
$ws.Range('C32').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -1.42%  '
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -0.49%  '
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +0.88%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.898'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -0.37%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.125.48'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -0.93%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.533'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -1.66%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.46'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -0.01%  '
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -1.07%  '
$ws.Range('B40').NumberFormat = '@'
$ws.Range('B40').Value = 'TrustWalletToken'
$ws.Range('C40').NumberFormat = '@'
$ws.Range('C40').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.794'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -0.48%  '
$ws.Range('B41').NumberFormat = '@'
$ws.Range('B41').Value = 'Quant'
$ws.Range('C41').NumberFormat = '@'
$ws.Range('C41').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '98.23'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -1.09%  '
$ws.Range('B42').NumberFormat = '@'
$ws.Range('B42').Value = 'FraxShare'
$ws.Range('C42').NumberFormat = '@'
$ws.Range('C42').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.32'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -2.75%  '
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -0.80%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '55.94'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -1.02%  '
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +2.72%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0520'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -2.17%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '7.72'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +1.13%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.413'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -0.32%  '
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +0.22%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0939'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -2.69%  '
$ws.Range('B51').NumberFormat = '@'
$ws.Range('B51').Value = 'NEARProtocol'
$ws.Range('C51').NumberFormat = '@'
$ws.Range('C51').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.16'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +0.75%  '
